# "combine two excel files into 1"
# - rename Sheet1 -> schedule
# - add a new sheet "desc" right after it, holding the teacher/class lookup
#   table that used to live in a second workbook
# - fix a typo in the schedule grid (DD1 -> Đ1)
# - restore the original selections on both sheets

$wb = $excel.ActiveWorkbook

# --- sheet1: rename to "schedule" -------------------------------------------------
$schedule = $wb.Worksheets.Item(1)
$schedule.Name = "schedule"

# --- add sheet2 "desc" right after "schedule" -------------------------------------
$desc = $wb.Worksheets.Add($null, $schedule)
$desc.Name = "desc"

# --- fix the one real data typo on the schedule sheet -----------------------------
$schedule.Range("B5").Value = "Đ1"

# --- populate the "desc" lookup table ---------------------------------------------
$header = @("TEACHER", "NAME", "CODE", "ID", "PASS")
for ($c = 0; $c -lt $header.Length; $c++) {
    $desc.Cells.Item(1, $c + 1).Value = $header[$c]
}

$rows = @(
    @("Triệu Lê Quang",        "Vật Lý",              "L10",  3141537349, 994494),
    @("Hạ Vũ Anh",             "Toán",                "T1",   8258073057, "Covid-19"),
    @("Khương Thị Thu Cúc",    "Ngữ Văn",              "V4",   9765088770, 614989),
    @("Nguyễn Mạnh Hà",        "Địa Lý",               "Đ1",   6822045197, 1234567890),
    @("Trần Văn Năng",         "Giáo dục công dân",    "G1",   8856018255, 123456789),
    @("Mai Thành Sơn",         "Tiếng Anh",            "N18",  2655252337, 91089),
    @("Nguyễn Văn Quảng",      "Công Nghệ",            "CN3",  7361311308, 2003),
    @("Nguyễn Thị Thu Cúc",    "Hóa Học",              "H12",  5841603699, "Hoahoc11"),
    @("Huỳnh Thị Ái Tâm",      "Sinh Học",             "Sh5",  7174188443, 66886868),
    @("Nguyễn Thu Hương",      "Lịch Sử",              "S5",   7666693318, "huong24"),
    @("Trần Mạnh Hùng",        "Thể dục",              "Td3",  5678955431, 989938),
    @("Nguyễn Văn Mạnh",       "Quốc Phòng",           "QP2",  9102597428, "thaymanhdz"),
    @("Bùi Tiến Dũng",         "Tin Học",              "Tin5", 2066767262, 123456)
)

$r = 2
foreach ($row in $rows) {
    $desc.Cells.Item($r, 1).Value = $row[0]
    $desc.Cells.Item($r, 2).Value = $row[1]
    $desc.Cells.Item($r, 3).Value = $row[2]
    $desc.Cells.Item($r, 4).Value = $row[3]
    $desc.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# NAME/CODE columns carry over the time-format style from the schedule sheet
$desc.Range("B1:C14").NumberFormat = "h:mm;@"

# --- selections -------------------------------------------------------------------
$desc.Range("A1:E14").Select() | Out-Null
$schedule.Range("J22").Select() | Out-Null
$schedule.Activate()
